$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (Volume number + reporting week dates) ---
$ws.Range("A8").Characters(21, 2).Text = "23"
$ws.Range("C9").Characters(27, 9).Text = "6/3/2024"
$ws.Range("C9").Characters(46, 8).Text = "6/9/2024"

# --- Crime-complaint table refresh (rows 15-31) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = 300
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I15").Value = 7
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 40
$ws.Range("L15").Value = 40

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = 16.129032258064
$ws.Range("L16").Value = 28.571428571428

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 91
$ws.Range("J17").Value = 103
$ws.Range("K17").Value = -11.650485436893
$ws.Range("L17").Value = -16.513761467889

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("E18").Value = 0
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 80
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = 40.90909090909
$ws.Range("L18").Value = -16.216216216216

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -47.619047619047
$ws.Range("I19").Value = 151
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = -11.176470588235
$ws.Range("L19").Value = -27.403846153846

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = -33.333333333333
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 114.285714285714
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = 47
$ws.Range("K20").Value = -29.787234042553
$ws.Range("L20").Value = -8.333333333333

# Row 21
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -36.363636363636
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -17.283950617283
$ws.Range("I21").Value = 349
$ws.Range("J21").Value = 378
$ws.Range("K21").Value = -7.671957671957
$ws.Range("L21").Value = -17.494089834515

# Row 23
$ws.Range("L23").Value = -76.923076923076

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -3.030303030303
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 123
$ws.Range("H24").Value = -6.50406504065
$ws.Range("I24").Value = 657
$ws.Range("J24").Value = 578
$ws.Range("K24").Value = 13.667820069204
$ws.Range("L24").Value = 15.263157894736

# Row 25
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -5.263157894736
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = 27.118644067796
$ws.Range("I25").Value = 450
$ws.Range("J25").Value = 346
$ws.Range("K25").Value = 30.057803468208
$ws.Range("L25").Value = 55.172413793103

# Row 26
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 266.666666666667
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -2.564102564102
$ws.Range("I26").Value = 262
$ws.Range("J26").Value = 238
$ws.Range("K26").Value = 10.084033613445
$ws.Range("L26").Value = 23.584905660377

# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 100
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("H27").Value = 800
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I27").Value = 15
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 114.285714285714

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 7
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 16.666666666666
$ws.Range("I28").Value = 27
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = -20.588235294117
$ws.Range("L28").Value = 22.727272727272

# Row 31
$ws.Range("C31").Value = 1
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = 0
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 5
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = 150
$ws.Range("L31").Value = 0

